$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.904.13"
$ws.Range("E2").Value = "  +2.03%  "
$ws.Range("D3").Value = "3.481.30"
$ws.Range("E3").Value = "  +2.53%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.02"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.61"
$ws.Range("E6").Value = "  +4.33%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.481"
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.66"
$ws.Range("E9").Value = "  -1.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.126"
$ws.Range("E10").Value = "  +2.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.401"
$ws.Range("E11").Value = "  +3.29%  "
$ws.Range("D12").Value = "4.076.78"
$ws.Range("E12").Value = "  +2.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "30.01"
$ws.Range("E13").Value = "  +5.36%  "
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("D15").Value = "3.483.33"
$ws.Range("E15").Value = "  +2.66%  "
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("D17").Value = "62.948.45"
$ws.Range("E17").Value = "  +2.13%  "
$ws.Range("E18").Value = "  +2.68%  "
$ws.Range("E19").Value = "  +5.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.35"
$ws.Range("E20").Value = "  +3.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "390.35"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.568"
$ws.Range("E22").Value = "  +2.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "75.19"
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("E25").Value = "  +2.60%  "
$ws.Range("E26").Value = "  +2.92%  "
$ws.Range("E27").Value = "  -7.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.69"
$ws.Range("E28").Value = "  +5.46%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.25"
$ws.Range("E30").Value = "  +2.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.14"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.41"
$ws.Range("E32").Value = "  +2.65%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.85"
$ws.Range("E34").Value = "  +1.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.13"
$ws.Range("E35").Value = "  +2.60%  "
$ws.Range("E36").Value = "  +3.87%  "
$ws.Range("E37").Value = "  +22.28%  "
$ws.Range("E38").Value = "  +2.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.58"
$ws.Range("E39").Value = "  +7.16%  "
$ws.Range("E40").Value = "  +2.66%  "
$ws.Range("E41").Value = "  -0.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.806"
$ws.Range("E42").Value = "  +3.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.28"
$ws.Range("E43").Value = "  -0.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.48"
$ws.Range("E44").Value = "  +0.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.72"
$ws.Range("E45").Value = "  +3.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.20"
$ws.Range("E46").Value = "  +2.99%  "
$ws.Range("D47").Value = "2.608.78"
$ws.Range("E47").Value = "  +5.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.63"
$ws.Range("E48").Value = "  +2.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.26"
$ws.Range("E49").Value = "  +9.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.81"
$ws.Range("E50").Value = "  +2.06%  "
$ws.Range("E51").Value = "  +0.09%  "
